$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" column header in F1, matching the style of the
# existing header row (copy format from E1, the last existing header cell).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the time_taken values for each data row (rows 2-8).
$ws.Range("F2").Value = "2021-10-05 13:39:09.324167"
$ws.Range("F3").Value = "2021-10-05 13:39:09.324178"
$ws.Range("F4").Value = "2021-10-05 13:39:09.324181"
$ws.Range("F5").Value = "2021-10-05 13:39:09.324184"
$ws.Range("F6").Value = "2021-10-05 13:39:09.324186"
$ws.Range("F7").Value = "2021-10-05 13:39:09.324189"
$ws.Range("F8").Value = "2021-10-05 13:39:09.324192"
